$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----
$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "Password"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Contact Number"
$ws.Range("E1").Value = "Address"
$ws.Range("F1").Value = "State"
$ws.Range("G1").Value = "Company Name"

# ---- Data row (row 2) ----
$ws.Range("A2").Value = "Test1"

# Numeric-looking text must stay text (not be coerced to a Number),
# so force a text format before assigning, then drop back to the
# Normal style so the cell keeps the default style index.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "123"
$ws.Range("B2").Style = "Normal"

$ws.Range("C2").Value = "TestUser"

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "123411313612"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = "blaabla"
$ws.Range("F2").Value = "Goa"
$ws.Range("G2").Value = "TestComapny"

# ---- Column widths ----
$ws.Columns.Item(1).ColumnWidth = 12
$ws.Columns.Item(2).ColumnWidth = 9.666666666666666
$ws.Columns.Item(4).ColumnWidth = 15.166666666666666
$ws.Columns.Item(5).ColumnWidth = 24
$ws.Columns.Item(6).ColumnWidth = 12.5
$ws.Columns.Item(7).ColumnWidth = 15.833333333333334

# ---- Selection ----
[void]$ws.Range("B2").Select()
